$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: ECs -> ECs
$ws.Cells.Item(2, 1).Value = "ECs"
$ws.Cells.Item(2, 2).Value = "C1qb"
$ws.Cells.Item(2, 3).Value = "Lrp1"
$ws.Cells.Item(2, 4).Value = "ECs"
$ws.Cells.Item(2, 5).Value = 1
$ws.Cells.Item(2, 6).Value = 0.5
$ws.Cells.Item(2, 7).Value = 0.7975575
$ws.Cells.Item(2, 8).Value = 1.595115
$ws.Cells.Item(2, 9).Value = 0.001587108162585159
$ws.Cells.Item(2, 10).Value = 0.001060423178847864
$ws.Cells.Item(2, 11).Value = 2
$ws.Cells.Item(2, 12).Value = 1
$ws.Cells.Item(2, 13).Value = 10.836393
$ws.Cells.Item(2, 14).Value = 21.672786
$ws.Cells.Item(2, 15).Value = 0.01846862257356514
$ws.Cells.Item(2, 16).Value = 0.01262296693390161
$ws.Cells.Item(2, 17).Value = 8.642646510097501
$ws.Cells.Item(2, 18).Value = 34.57058604039
$ws.Cells.Item(2, 19).Value = 0.00002931170163820975
$ws.Cells.Item(2, 20).Value = 0.00001338568672253942

# Row 3: ECs -> FAPs
$ws.Cells.Item(3, 1).Value = "ECs"
$ws.Cells.Item(3, 2).Value = "C1qb"
$ws.Cells.Item(3, 3).Value = "Lrp1"
$ws.Cells.Item(3, 4).Value = "FAPs"
$ws.Cells.Item(3, 5).Value = 1
$ws.Cells.Item(3, 6).Value = 0.5
$ws.Cells.Item(3, 7).Value = 0.7975575
$ws.Cells.Item(3, 8).Value = 1.595115
$ws.Cells.Item(3, 9).Value = 0.001587108162585159
$ws.Cells.Item(3, 10).Value = 0.001060423178847864
$ws.Cells.Item(3, 11).Value = 3
$ws.Cells.Item(3, 12).Value = 1
$ws.Cells.Item(3, 13).Value = 139.6948166666666
$ws.Cells.Item(3, 14).Value = 419.0844499999999
$ws.Cells.Item(3, 15).Value = 0.2380839126543345
$ws.Cells.Item(3, 16).Value = 0.2440890227431923
$ws.Cells.Item(3, 17).Value = 111.414648743625
$ws.Cells.Item(3, 18).Value = 668.48789246175
$ws.Cells.Item(3, 19).Value = 0.0003778649211539063
$ws.Cells.Item(3, 20).Value = 0.0002588376574192046

# Row 4: ECs -> Inflammatory-Mac
$ws.Cells.Item(4, 1).Value = "ECs"
$ws.Cells.Item(4, 2).Value = "C1qb"
$ws.Cells.Item(4, 3).Value = "Lrp1"
$ws.Cells.Item(4, 4).Value = "Inflammatory-Mac"
$ws.Cells.Item(4, 5).Value = 1
$ws.Cells.Item(4, 6).Value = 0.5
$ws.Cells.Item(4, 7).Value = 0.7975575
$ws.Cells.Item(4, 8).Value = 1.595115
$ws.Cells.Item(4, 9).Value = 0.001587108162585159
$ws.Cells.Item(4, 10).Value = 0.001060423178847864
$ws.Cells.Item(4, 11).Value = 3
$ws.Cells.Item(4, 12).Value = 1
$ws.Cells.Item(4, 13).Value = 150.12088
$ws.Cells.Item(4, 14).Value = 450.36264
$ws.Cells.Item(4, 15).Value = 0.2558532043948076
$ws.Cells.Item(4, 16).Value = 0.2623065033256284
$ws.Cells.Item(4, 17).Value = 119.7300337506
$ws.Cells.Item(4, 18).Value = 718.3802025036
$ws.Cells.Item(4, 19).Value = 0.0004060667091185682
$ws.Cells.Item(4, 20).Value = 0.0002781558960890308

# Row 5: ECs -> MuSCs
$ws.Cells.Item(5, 1).Value = "ECs"
$ws.Cells.Item(5, 2).Value = "C1qb"
$ws.Cells.Item(5, 3).Value = "Lrp1"
$ws.Cells.Item(5, 4).Value = "MuSCs"
$ws.Cells.Item(5, 5).Value = 1
$ws.Cells.Item(5, 6).Value = 0.5
$ws.Cells.Item(5, 7).Value = 0.7975575
$ws.Cells.Item(5, 8).Value = 1.595115
$ws.Cells.Item(5, 9).Value = 0.001587108162585159
$ws.Cells.Item(5, 10).Value = 0.001060423178847864
$ws.Cells.Item(5, 11).Value = 2
$ws.Cells.Item(5, 12).Value = 1
$ws.Cells.Item(5, 13).Value = 32.469223
$ws.Cells.Item(5, 14).Value = 64.938446
$ws.Cells.Item(5, 15).Value = 0.05533777012737728
$ws.Cells.Item(5, 16).Value = 0.03782235733776705
$ws.Cells.Item(5, 17).Value = 25.8960723228225
$ws.Cells.Item(5, 18).Value = 103.58428929129
$ws.Cells.Item(5, 19).Value = 0.00008782702666842163
$ws.Cells.Item(5, 20).Value = 0.00004010770439963478

# Row 6: ECs -> Neutrophils
$ws.Cells.Item(6, 1).Value = "ECs"
$ws.Cells.Item(6, 2).Value = "C1qb"
$ws.Cells.Item(6, 3).Value = "Lrp1"
$ws.Cells.Item(6, 4).Value = "Neutrophils"
$ws.Cells.Item(6, 5).Value = 1
$ws.Cells.Item(6, 6).Value = 0.5
$ws.Cells.Item(6, 7).Value = 0.7975575
$ws.Cells.Item(6, 8).Value = 1.595115
$ws.Cells.Item(6, 9).Value = 0.001587108162585159
$ws.Cells.Item(6, 10).Value = 0.001060423178847864
$ws.Cells.Item(6, 11).Value = 3
$ws.Cells.Item(6, 12).Value = 1
$ws.Cells.Item(6, 13).Value = 52.09024033333333
$ws.Cells.Item(6, 14).Value = 156.270721
$ws.Cells.Item(6, 15).Value = 0.08877815602319267
$ws.Cells.Item(6, 16).Value = 0.09101737745760805
$ws.Cells.Item(6, 17).Value = 41.5449618546525
$ws.Cells.Item(6, 18).Value = 249.269771127915
$ws.Cells.Item(6, 19).Value = 0.0001409005360836679
$ws.Cells.Item(6, 20).Value = 0.00009651693673399267

# Row 7: ECs -> Resolving-Mac
$ws.Cells.Item(7, 1).Value = "ECs"
$ws.Cells.Item(7, 2).Value = "C1qb"
$ws.Cells.Item(7, 3).Value = "Lrp1"
$ws.Cells.Item(7, 4).Value = "Resolving-Mac"
$ws.Cells.Item(7, 5).Value = 1
$ws.Cells.Item(7, 6).Value = 0.5
$ws.Cells.Item(7, 7).Value = 0.7975575
$ws.Cells.Item(7, 8).Value = 1.595115
$ws.Cells.Item(7, 9).Value = 0.001587108162585159
$ws.Cells.Item(7, 10).Value = 0.001060423178847864
$ws.Cells.Item(7, 11).Value = 3
$ws.Cells.Item(7, 12).Value = 1
$ws.Cells.Item(7, 13).Value = 201.5345866666667
$ws.Cells.Item(7, 14).Value = 604.60376
$ws.Cells.Item(7, 15).Value = 0.3434783342267227
$ws.Cells.Item(7, 16).Value = 0.3521417722019025
$ws.Cells.Item(7, 17).Value = 160.7354211054
$ws.Cells.Item(7, 18).Value = 964.4125266323999
$ws.Cells.Item(7, 19).Value = 0.000545137267922385
$ws.Cells.Item(7, 20).Value = 0.000373419297483462

# Row 8: Inflammatory-Mac -> ECs
$ws.Cells.Item(8, 1).Value = "Inflammatory-Mac"
$ws.Cells.Item(8, 2).Value = "C1qb"
$ws.Cells.Item(8, 3).Value = "Lrp1"
$ws.Cells.Item(8, 4).Value = "ECs"
$ws.Cells.Item(8, 5).Value = 3
$ws.Cells.Item(8, 6).Value = 1
$ws.Cells.Item(8, 7).Value = 170.2842533333333
$ws.Cells.Item(8, 8).Value = 510.85276
$ws.Cells.Item(8, 9).Value = 0.3388589893832762
$ws.Cells.Item(8, 10).Value = 0.3396119450211458
$ws.Cells.Item(8, 11).Value = 2
$ws.Cells.Item(8, 12).Value = 1
$ws.Cells.Item(8, 13).Value = 10.836393
$ws.Cells.Item(8, 14).Value = 21.672786
$ws.Cells.Item(8, 15).Value = 0.01846862257356514
$ws.Cells.Item(8, 16).Value = 0.01262296693390161
$ws.Cells.Item(8, 17).Value = 1845.26709083156
$ws.Cells.Item(8, 18).Value = 11071.60254498936
$ws.Cells.Item(8, 19).Value = 0.006258258780579445
$ws.Cells.Item(8, 20).Value = 0.004286910352359935

# Row 9: Inflammatory-Mac -> FAPs
$ws.Cells.Item(9, 1).Value = "Inflammatory-Mac"
$ws.Cells.Item(9, 2).Value = "C1qb"
$ws.Cells.Item(9, 3).Value = "Lrp1"
$ws.Cells.Item(9, 4).Value = "FAPs"
$ws.Cells.Item(9, 5).Value = 3
$ws.Cells.Item(9, 6).Value = 1
$ws.Cells.Item(9, 7).Value = 170.2842533333333
$ws.Cells.Item(9, 8).Value = 510.85276
$ws.Cells.Item(9, 9).Value = 0.3388589893832762
$ws.Cells.Item(9, 10).Value = 0.3396119450211458
$ws.Cells.Item(9, 11).Value = 3
$ws.Cells.Item(9, 12).Value = 1
$ws.Cells.Item(9, 13).Value = 139.6948166666666
$ws.Cells.Item(9, 14).Value = 419.0844499999999
$ws.Cells.Item(9, 15).Value = 0.2380839126543345
$ws.Cells.Item(9, 16).Value = 0.2440890227431923
$ws.Cells.Item(9, 17).Value = 23787.82755062022
$ws.Cells.Item(9, 18).Value = 214090.447955582
$ws.Cells.Item(9, 19).Value = 0.08067687403046402
$ws.Cells.Item(9, 20).Value = 0.08289554777212625

# Row 10: Inflammatory-Mac -> Inflammatory-Mac
$ws.Cells.Item(10, 1).Value = "Inflammatory-Mac"
$ws.Cells.Item(10, 2).Value = "C1qb"
$ws.Cells.Item(10, 3).Value = "Lrp1"
$ws.Cells.Item(10, 4).Value = "Inflammatory-Mac"
$ws.Cells.Item(10, 5).Value = 3
$ws.Cells.Item(10, 6).Value = 1
$ws.Cells.Item(10, 7).Value = 170.2842533333333
$ws.Cells.Item(10, 8).Value = 510.85276
$ws.Cells.Item(10, 9).Value = 0.3388589893832762
$ws.Cells.Item(10, 10).Value = 0.3396119450211458
$ws.Cells.Item(10, 11).Value = 3
$ws.Cells.Item(10, 12).Value = 1
$ws.Cells.Item(10, 13).Value = 150.12088
$ws.Cells.Item(10, 14).Value = 450.36264
$ws.Cells.Item(10, 15).Value = 0.2558532043948076
$ws.Cells.Item(10, 16).Value = 0.2623065033256284
$ws.Cells.Item(10, 17).Value = 25563.22196054293
$ws.Cells.Item(10, 18).Value = 230068.9976448864
$ws.Cells.Item(10, 19).Value = 0.08669815827169733
$ws.Cells.Item(10, 20).Value = 0.08908242178611234

# Row 11: Inflammatory-Mac -> MuSCs
$ws.Cells.Item(11, 1).Value = "Inflammatory-Mac"
$ws.Cells.Item(11, 2).Value = "C1qb"
$ws.Cells.Item(11, 3).Value = "Lrp1"
$ws.Cells.Item(11, 4).Value = "MuSCs"
$ws.Cells.Item(11, 5).Value = 3
$ws.Cells.Item(11, 6).Value = 1
$ws.Cells.Item(11, 7).Value = 170.2842533333333
$ws.Cells.Item(11, 8).Value = 510.85276
$ws.Cells.Item(11, 9).Value = 0.3388589893832762
$ws.Cells.Item(11, 10).Value = 0.3396119450211458
$ws.Cells.Item(11, 11).Value = 2
$ws.Cells.Item(11, 12).Value = 1
$ws.Cells.Item(11, 13).Value = 32.469223
$ws.Cells.Item(11, 14).Value = 64.938446
$ws.Cells.Item(11, 15).Value = 0.05533777012737728
$ws.Cells.Item(11, 16).Value = 0.03782235733776705
$ws.Cells.Item(11, 17).Value = 5528.997394868494
$ws.Cells.Item(11, 18).Value = 33173.98436921096
$ws.Cells.Item(11, 19).Value = 0.01875170086008712
$ws.Cells.Item(11, 20).Value = 0.01284492434076388

# Row 12: Inflammatory-Mac -> Neutrophils
$ws.Cells.Item(12, 1).Value = "Inflammatory-Mac"
$ws.Cells.Item(12, 2).Value = "C1qb"
$ws.Cells.Item(12, 3).Value = "Lrp1"
$ws.Cells.Item(12, 4).Value = "Neutrophils"
$ws.Cells.Item(12, 5).Value = 3
$ws.Cells.Item(12, 6).Value = 1
$ws.Cells.Item(12, 7).Value = 170.2842533333333
$ws.Cells.Item(12, 8).Value = 510.85276
$ws.Cells.Item(12, 9).Value = 0.3388589893832762
$ws.Cells.Item(12, 10).Value = 0.3396119450211458
$ws.Cells.Item(12, 11).Value = 3
$ws.Cells.Item(12, 12).Value = 1
$ws.Cells.Item(12, 13).Value = 52.09024033333333
$ws.Cells.Item(12, 14).Value = 156.270721
$ws.Cells.Item(12, 15).Value = 0.08877815602319267
$ws.Cells.Item(12, 16).Value = 0.09101737745760805
$ws.Cells.Item(12, 17).Value = 8870.14768111555
$ws.Cells.Item(12, 18).Value = 79831.32913003994
$ws.Cells.Item(12, 19).Value = 0.03008327622932989
$ws.Cells.Item(12, 20).Value = 0.03091058858910206

# Row 13: Inflammatory-Mac -> Resolving-Mac
$ws.Cells.Item(13, 1).Value = "Inflammatory-Mac"
$ws.Cells.Item(13, 2).Value = "C1qb"
$ws.Cells.Item(13, 3).Value = "Lrp1"
$ws.Cells.Item(13, 4).Value = "Resolving-Mac"
$ws.Cells.Item(13, 5).Value = 3
$ws.Cells.Item(13, 6).Value = 1
$ws.Cells.Item(13, 7).Value = 170.2842533333333
$ws.Cells.Item(13, 8).Value = 510.85276
$ws.Cells.Item(13, 9).Value = 0.3388589893832762
$ws.Cells.Item(13, 10).Value = 0.3396119450211458
$ws.Cells.Item(13, 11).Value = 3
$ws.Cells.Item(13, 12).Value = 1
$ws.Cells.Item(13, 13).Value = 201.5345866666667
$ws.Cells.Item(13, 14).Value = 604.60376
$ws.Cells.Item(13, 15).Value = 0.3434783342267227
$ws.Cells.Item(13, 16).Value = 0.3521417722019025
$ws.Cells.Item(13, 17).Value = 34318.16661137529
$ws.Cells.Item(13, 18).Value = 308863.4995023776
$ws.Cells.Item(13, 19).Value = 0.1163907212111185
$ws.Cells.Item(13, 20).Value = 0.1195915521806814

# Row 14: MuSCs -> ECs
$ws.Cells.Item(14, 1).Value = "MuSCs"
$ws.Cells.Item(14, 2).Value = "C1qb"
$ws.Cells.Item(14, 3).Value = "Lrp1"
$ws.Cells.Item(14, 4).Value = "ECs"
$ws.Cells.Item(14, 5).Value = 1
$ws.Cells.Item(14, 6).Value = 0.5
$ws.Cells.Item(14, 7).Value = 2.544879
$ws.Cells.Item(14, 8).Value = 5.089758
$ws.Cells.Item(14, 9).Value = 0.005064209456611662
$ws.Cells.Item(14, 10).Value = 0.003383641529247952
$ws.Cells.Item(14, 11).Value = 2
$ws.Cells.Item(14, 12).Value = 1
$ws.Cells.Item(14, 13).Value = 10.836393
$ws.Cells.Item(14, 14).Value = 21.672786
$ws.Cells.Item(14, 15).Value = 0.01846862257356514
$ws.Cells.Item(14, 16).Value = 0.01262296693390161
$ws.Cells.Item(14, 17).Value = 27.577308981447
$ws.Cells.Item(14, 18).Value = 110.309235925788
$ws.Cells.Item(14, 19).Value = 0.00009352897308764017
$ws.Cells.Item(14, 20).Value = 0.00004271159513987317

# Row 15: MuSCs -> FAPs
$ws.Cells.Item(15, 1).Value = "MuSCs"
$ws.Cells.Item(15, 2).Value = "C1qb"
$ws.Cells.Item(15, 3).Value = "Lrp1"
$ws.Cells.Item(15, 4).Value = "FAPs"
$ws.Cells.Item(15, 5).Value = 1
$ws.Cells.Item(15, 6).Value = 0.5
$ws.Cells.Item(15, 7).Value = 2.544879
$ws.Cells.Item(15, 8).Value = 5.089758
$ws.Cells.Item(15, 9).Value = 0.005064209456611662
$ws.Cells.Item(15, 10).Value = 0.003383641529247952
$ws.Cells.Item(15, 11).Value = 3
$ws.Cells.Item(15, 12).Value = 1
$ws.Cells.Item(15, 13).Value = 139.6948166666666
$ws.Cells.Item(15, 14).Value = 419.0844499999999
$ws.Cells.Item(15, 15).Value = 0.2380839126543345
$ws.Cells.Item(15, 16).Value = 0.2440890227431923
$ws.Cells.Item(15, 17).Value = 355.5064053438499
$ws.Cells.Item(15, 18).Value = 2133.0384320631
$ws.Cells.Item(15, 19).Value = 0.001205706801931186
$ws.Cells.Item(15, 20).Value = 0.0008259097541874136

# Row 16: MuSCs -> Inflammatory-Mac
$ws.Cells.Item(16, 1).Value = "MuSCs"
$ws.Cells.Item(16, 2).Value = "C1qb"
$ws.Cells.Item(16, 3).Value = "Lrp1"
$ws.Cells.Item(16, 4).Value = "Inflammatory-Mac"
$ws.Cells.Item(16, 5).Value = 1
$ws.Cells.Item(16, 6).Value = 0.5
$ws.Cells.Item(16, 7).Value = 2.544879
$ws.Cells.Item(16, 8).Value = 5.089758
$ws.Cells.Item(16, 9).Value = 0.005064209456611662
$ws.Cells.Item(16, 10).Value = 0.003383641529247952
$ws.Cells.Item(16, 11).Value = 3
$ws.Cells.Item(16, 12).Value = 1
$ws.Cells.Item(16, 13).Value = 150.12088
$ws.Cells.Item(16, 14).Value = 450.36264
$ws.Cells.Item(16, 15).Value = 0.2558532043948076
$ws.Cells.Item(16, 16).Value = 0.2623065033256284
$ws.Cells.Item(16, 17).Value = 382.03947497352
$ws.Cells.Item(16, 18).Value = 2292.23684984112
$ws.Cells.Item(16, 19).Value = 0.001295694217200581
$ws.Cells.Item(16, 20).Value = 0.0008875511780444125

# Row 17: MuSCs -> MuSCs
$ws.Cells.Item(17, 1).Value = "MuSCs"
$ws.Cells.Item(17, 2).Value = "C1qb"
$ws.Cells.Item(17, 3).Value = "Lrp1"
$ws.Cells.Item(17, 4).Value = "MuSCs"
$ws.Cells.Item(17, 5).Value = 1
$ws.Cells.Item(17, 6).Value = 0.5
$ws.Cells.Item(17, 7).Value = 2.544879
$ws.Cells.Item(17, 8).Value = 5.089758
$ws.Cells.Item(17, 9).Value = 0.005064209456611662
$ws.Cells.Item(17, 10).Value = 0.003383641529247952
$ws.Cells.Item(17, 11).Value = 2
$ws.Cells.Item(17, 12).Value = 1
$ws.Cells.Item(17, 13).Value = 32.469223
$ws.Cells.Item(17, 14).Value = 64.938446
$ws.Cells.Item(17, 15).Value = 0.05533777012737728
$ws.Cells.Item(17, 16).Value = 0.03782235733776705
$ws.Cells.Item(17, 17).Value = 82.63024375901699
$ws.Cells.Item(17, 18).Value = 330.520975036068
$ws.Cells.Item(17, 19).Value = 0.0002802420587868663
$ws.Cells.Item(17, 20).Value = 0.0001279772990221246

# Row 18: MuSCs -> Neutrophils
$ws.Cells.Item(18, 1).Value = "MuSCs"
$ws.Cells.Item(18, 2).Value = "C1qb"
$ws.Cells.Item(18, 3).Value = "Lrp1"
$ws.Cells.Item(18, 4).Value = "Neutrophils"
$ws.Cells.Item(18, 5).Value = 1
$ws.Cells.Item(18, 6).Value = 0.5
$ws.Cells.Item(18, 7).Value = 2.544879
$ws.Cells.Item(18, 8).Value = 5.089758
$ws.Cells.Item(18, 9).Value = 0.005064209456611662
$ws.Cells.Item(18, 10).Value = 0.003383641529247952
$ws.Cells.Item(18, 11).Value = 3
$ws.Cells.Item(18, 12).Value = 1
$ws.Cells.Item(18, 13).Value = 52.09024033333333
$ws.Cells.Item(18, 14).Value = 156.270721
$ws.Cells.Item(18, 15).Value = 0.08877815602319267
$ws.Cells.Item(18, 16).Value = 0.09101737745760805
$ws.Cells.Item(18, 17).Value = 132.563358729253
$ws.Cells.Item(18, 18).Value = 795.3801523755179
$ws.Cells.Item(18, 19).Value = 0.0004495911772731979
$ws.Cells.Item(18, 20).Value = 0.000307970178248799

# Row 19: MuSCs -> Resolving-Mac
$ws.Cells.Item(19, 1).Value = "MuSCs"
$ws.Cells.Item(19, 2).Value = "C1qb"
$ws.Cells.Item(19, 3).Value = "Lrp1"
$ws.Cells.Item(19, 4).Value = "Resolving-Mac"
$ws.Cells.Item(19, 5).Value = 1
$ws.Cells.Item(19, 6).Value = 0.5
$ws.Cells.Item(19, 7).Value = 2.544879
$ws.Cells.Item(19, 8).Value = 5.089758
$ws.Cells.Item(19, 9).Value = 0.005064209456611662
$ws.Cells.Item(19, 10).Value = 0.003383641529247952
$ws.Cells.Item(19, 11).Value = 3
$ws.Cells.Item(19, 12).Value = 1
$ws.Cells.Item(19, 13).Value = 201.5345866666667
$ws.Cells.Item(19, 14).Value = 604.60376
$ws.Cells.Item(19, 15).Value = 0.3434783342267227
$ws.Cells.Item(19, 16).Value = 0.3521417722019025
$ws.Cells.Item(19, 17).Value = 512.88113738168
$ws.Cells.Item(19, 18).Value = 3077.286824290079
$ws.Cells.Item(19, 19).Value = 0.00173944622833219
$ws.Cells.Item(19, 20).Value = 0.00119152152460533

# Row 20: Neutrophils -> ECs
$ws.Cells.Item(20, 1).Value = "Neutrophils"
$ws.Cells.Item(20, 2).Value = "C1qb"
$ws.Cells.Item(20, 3).Value = "Lrp1"
$ws.Cells.Item(20, 4).Value = "ECs"
$ws.Cells.Item(20, 5).Value = 3
$ws.Cells.Item(20, 6).Value = 1
$ws.Cells.Item(20, 7).Value = 51.68882633333334
$ws.Cells.Item(20, 8).Value = 155.066479
$ws.Cells.Item(20, 9).Value = 0.1028587383205349
$ws.Cells.Item(20, 10).Value = 0.1030872937649797
$ws.Cells.Item(20, 11).Value = 2
$ws.Cells.Item(20, 12).Value = 1
$ws.Cells.Item(20, 13).Value = 10.836393
$ws.Cells.Item(20, 14).Value = 21.672786
$ws.Cells.Item(20, 15).Value = 0.01846862257356514
$ws.Cells.Item(20, 16).Value = 0.01262296693390161
$ws.Cells.Item(20, 17).Value = 560.120435856749
$ws.Cells.Item(20, 18).Value = 3360.722615140495
$ws.Cells.Item(20, 19).Value = 0.001899659216435061
$ws.Cells.Item(20, 20).Value = 0.00130126750050074

# Row 21: Neutrophils -> FAPs
$ws.Cells.Item(21, 1).Value = "Neutrophils"
$ws.Cells.Item(21, 2).Value = "C1qb"
$ws.Cells.Item(21, 3).Value = "Lrp1"
$ws.Cells.Item(21, 4).Value = "FAPs"
$ws.Cells.Item(21, 5).Value = 3
$ws.Cells.Item(21, 6).Value = 1
$ws.Cells.Item(21, 7).Value = 51.68882633333334
$ws.Cells.Item(21, 8).Value = 155.066479
$ws.Cells.Item(21, 9).Value = 0.1028587383205349
$ws.Cells.Item(21, 10).Value = 0.1030872937649797
$ws.Cells.Item(21, 11).Value = 3
$ws.Cells.Item(21, 12).Value = 1
$ws.Cells.Item(21, 13).Value = 139.6948166666666
$ws.Cells.Item(21, 14).Value = 419.0844499999999
$ws.Cells.Item(21, 15).Value = 0.2380839126543345
$ws.Cells.Item(21, 16).Value = 0.2440890227431923
$ws.Cells.Item(21, 17).Value = 7220.661118350172
$ws.Cells.Item(21, 18).Value = 64985.95006515155
$ws.Cells.Item(21, 19).Value = 0.02448901087004129
$ws.Cells.Item(21, 20).Value = 0.02516247679233428

# Row 22: Neutrophils -> Inflammatory-Mac
$ws.Cells.Item(22, 1).Value = "Neutrophils"
$ws.Cells.Item(22, 2).Value = "C1qb"
$ws.Cells.Item(22, 3).Value = "Lrp1"
$ws.Cells.Item(22, 4).Value = "Inflammatory-Mac"
$ws.Cells.Item(22, 5).Value = 3
$ws.Cells.Item(22, 6).Value = 1
$ws.Cells.Item(22, 7).Value = 51.68882633333334
$ws.Cells.Item(22, 8).Value = 155.066479
$ws.Cells.Item(22, 9).Value = 0.1028587383205349
$ws.Cells.Item(22, 10).Value = 0.1030872937649797
$ws.Cells.Item(22, 11).Value = 3
$ws.Cells.Item(22, 12).Value = 1
$ws.Cells.Item(22, 13).Value = 150.12088
$ws.Cells.Item(22, 14).Value = 450.36264
$ws.Cells.Item(22, 15).Value = 0.2558532043948076
$ws.Cells.Item(22, 16).Value = 0.2623065033256284
$ws.Cells.Item(22, 17).Value = 7759.572095327174
$ws.Cells.Item(22, 18).Value = 69836.14885794457
$ws.Cells.Item(22, 19).Value = 0.02631673779931585
$ws.Cells.Item(22, 20).Value = 0.02704046756479368

# Row 23: Neutrophils -> MuSCs
$ws.Cells.Item(23, 1).Value = "Neutrophils"
$ws.Cells.Item(23, 2).Value = "C1qb"
$ws.Cells.Item(23, 3).Value = "Lrp1"
$ws.Cells.Item(23, 4).Value = "MuSCs"
$ws.Cells.Item(23, 5).Value = 3
$ws.Cells.Item(23, 6).Value = 1
$ws.Cells.Item(23, 7).Value = 51.68882633333334
$ws.Cells.Item(23, 8).Value = 155.066479
$ws.Cells.Item(23, 9).Value = 0.1028587383205349
$ws.Cells.Item(23, 10).Value = 0.1030872937649797
$ws.Cells.Item(23, 11).Value = 2
$ws.Cells.Item(23, 12).Value = 1
$ws.Cells.Item(23, 13).Value = 32.469223
$ws.Cells.Item(23, 14).Value = 64.938446
$ws.Cells.Item(23, 15).Value = 0.05533777012737728
$ws.Cells.Item(23, 16).Value = 0.03782235733776705
$ws.Cells.Item(23, 17).Value = 1678.296028825273
$ws.Cells.Item(23, 18).Value = 10069.77617295164
$ws.Cells.Item(23, 19).Value = 0.005691973216773814
$ws.Cells.Item(23, 20).Value = 0.003899004461762428

# Row 24: Neutrophils -> Neutrophils
$ws.Cells.Item(24, 1).Value = "Neutrophils"
$ws.Cells.Item(24, 2).Value = "C1qb"
$ws.Cells.Item(24, 3).Value = "Lrp1"
$ws.Cells.Item(24, 4).Value = "Neutrophils"
$ws.Cells.Item(24, 5).Value = 3
$ws.Cells.Item(24, 6).Value = 1
$ws.Cells.Item(24, 7).Value = 51.68882633333334
$ws.Cells.Item(24, 8).Value = 155.066479
$ws.Cells.Item(24, 9).Value = 0.1028587383205349
$ws.Cells.Item(24, 10).Value = 0.1030872937649797
$ws.Cells.Item(24, 11).Value = 3
$ws.Cells.Item(24, 12).Value = 1
$ws.Cells.Item(24, 13).Value = 52.09024033333333
$ws.Cells.Item(24, 14).Value = 156.270721
$ws.Cells.Item(24, 15).Value = 0.08877815602319267
$ws.Cells.Item(24, 16).Value = 0.09101737745760805
$ws.Cells.Item(24, 17).Value = 2692.483386251262
$ws.Cells.Item(24, 18).Value = 24232.35047626136
$ws.Cells.Item(24, 19).Value = 0.009131609118969196
$ws.Cells.Item(24, 20).Value = 0.009382735127690481

# Row 25: Neutrophils -> Resolving-Mac
$ws.Cells.Item(25, 1).Value = "Neutrophils"
$ws.Cells.Item(25, 2).Value = "C1qb"
$ws.Cells.Item(25, 3).Value = "Lrp1"
$ws.Cells.Item(25, 4).Value = "Resolving-Mac"
$ws.Cells.Item(25, 5).Value = 3
$ws.Cells.Item(25, 6).Value = 1
$ws.Cells.Item(25, 7).Value = 51.68882633333334
$ws.Cells.Item(25, 8).Value = 155.066479
$ws.Cells.Item(25, 9).Value = 0.1028587383205349
$ws.Cells.Item(25, 10).Value = 0.1030872937649797
$ws.Cells.Item(25, 11).Value = 3
$ws.Cells.Item(25, 12).Value = 1
$ws.Cells.Item(25, 13).Value = 201.5345866666667
$ws.Cells.Item(25, 14).Value = 604.60376
$ws.Cells.Item(25, 15).Value = 0.3434783342267227
$ws.Cells.Item(25, 16).Value = 0.3521417722019025
$ws.Cells.Item(25, 17).Value = 10417.08625037345
$ws.Cells.Item(25, 18).Value = 93753.77625336104
$ws.Cells.Item(25, 19).Value = 0.03532974809899971
$ws.Cells.Item(25, 20).Value = 0.03630134231789809

# Row 26: Resolving-Mac -> ECs
$ws.Cells.Item(26, 1).Value = "Resolving-Mac"
$ws.Cells.Item(26, 2).Value = "C1qb"
$ws.Cells.Item(26, 3).Value = "Lrp1"
$ws.Cells.Item(26, 4).Value = "ECs"
$ws.Cells.Item(26, 5).Value = 3
$ws.Cells.Item(26, 6).Value = 1
$ws.Cells.Item(26, 7).Value = 277.206945
$ws.Cells.Item(26, 8).Value = 831.6208349999999
$ws.Cells.Item(26, 9).Value = 0.551630954676992
$ws.Cells.Item(26, 10).Value = 0.5528566965057786
$ws.Cells.Item(26, 11).Value = 2
$ws.Cells.Item(26, 12).Value = 1
$ws.Cells.Item(26, 13).Value = 10.836393
$ws.Cells.Item(26, 14).Value = 21.672786
$ws.Cells.Item(26, 15).Value = 0.01846862257356514
$ws.Cells.Item(26, 16).Value = 0.01262296693390161
$ws.Cells.Item(26, 17).Value = 3003.923398349385
$ws.Cells.Item(26, 18).Value = 18023.54039009631
$ws.Cells.Item(26, 19).Value = 0.01018786390182478
$ws.Cells.Item(26, 20).Value = 0.006978691799178521

# Row 27: Resolving-Mac -> FAPs
$ws.Cells.Item(27, 1).Value = "Resolving-Mac"
$ws.Cells.Item(27, 2).Value = "C1qb"
$ws.Cells.Item(27, 3).Value = "Lrp1"
$ws.Cells.Item(27, 4).Value = "FAPs"
$ws.Cells.Item(27, 5).Value = 3
$ws.Cells.Item(27, 6).Value = 1
$ws.Cells.Item(27, 7).Value = 277.206945
$ws.Cells.Item(27, 8).Value = 831.6208349999999
$ws.Cells.Item(27, 9).Value = 0.551630954676992
$ws.Cells.Item(27, 10).Value = 0.5528566965057786
$ws.Cells.Item(27, 11).Value = 3
$ws.Cells.Item(27, 12).Value = 1
$ws.Cells.Item(27, 13).Value = 139.6948166666666
$ws.Cells.Item(27, 14).Value = 419.0844499999999
$ws.Cells.Item(27, 15).Value = 0.2380839126543345
$ws.Cells.Item(27, 16).Value = 0.2440890227431923
$ws.Cells.Item(27, 17).Value = 38724.37336050174
$ws.Cells.Item(27, 18).Value = 348519.3602445157
$ws.Cells.Item(27, 19).Value = 0.1313344560307441
$ws.Cells.Item(27, 20).Value = 0.1349462507671252

# Row 28: Resolving-Mac -> Inflammatory-Mac
$ws.Cells.Item(28, 1).Value = "Resolving-Mac"
$ws.Cells.Item(28, 2).Value = "C1qb"
$ws.Cells.Item(28, 3).Value = "Lrp1"
$ws.Cells.Item(28, 4).Value = "Inflammatory-Mac"
$ws.Cells.Item(28, 5).Value = 3
$ws.Cells.Item(28, 6).Value = 1
$ws.Cells.Item(28, 7).Value = 277.206945
$ws.Cells.Item(28, 8).Value = 831.6208349999999
$ws.Cells.Item(28, 9).Value = 0.551630954676992
$ws.Cells.Item(28, 10).Value = 0.5528566965057786
$ws.Cells.Item(28, 11).Value = 3
$ws.Cells.Item(28, 12).Value = 1
$ws.Cells.Item(28, 13).Value = 150.12088
$ws.Cells.Item(28, 14).Value = 450.36264
$ws.Cells.Item(28, 15).Value = 0.2558532043948076
$ws.Cells.Item(28, 16).Value = 0.2623065033256284
$ws.Cells.Item(28, 17).Value = 41614.55052551159
$ws.Cells.Item(28, 18).Value = 374530.9547296044
$ws.Cells.Item(28, 19).Value = 0.1411365473974753
$ws.Cells.Item(28, 20).Value = 0.145017906900589

# Row 29: Resolving-Mac -> MuSCs
$ws.Cells.Item(29, 1).Value = "Resolving-Mac"
$ws.Cells.Item(29, 2).Value = "C1qb"
$ws.Cells.Item(29, 3).Value = "Lrp1"
$ws.Cells.Item(29, 4).Value = "MuSCs"
$ws.Cells.Item(29, 5).Value = 3
$ws.Cells.Item(29, 6).Value = 1
$ws.Cells.Item(29, 7).Value = 277.206945
$ws.Cells.Item(29, 8).Value = 831.6208349999999
$ws.Cells.Item(29, 9).Value = 0.551630954676992
$ws.Cells.Item(29, 10).Value = 0.5528566965057786
$ws.Cells.Item(29, 11).Value = 2
$ws.Cells.Item(29, 12).Value = 1
$ws.Cells.Item(29, 13).Value = 32.469223
$ws.Cells.Item(29, 14).Value = 64.938446
$ws.Cells.Item(29, 15).Value = 0.05533777012737728
$ws.Cells.Item(29, 16).Value = 0.03782235733776705
$ws.Cells.Item(29, 17).Value = 9000.694114353733
$ws.Cells.Item(29, 18).Value = 54004.16468612241
$ws.Cells.Item(29, 19).Value = 0.03052602696506105
$ws.Cells.Item(29, 20).Value = 0.02091034353181899

# Row 30: Resolving-Mac -> Neutrophils
$ws.Cells.Item(30, 1).Value = "Resolving-Mac"
$ws.Cells.Item(30, 2).Value = "C1qb"
$ws.Cells.Item(30, 3).Value = "Lrp1"
$ws.Cells.Item(30, 4).Value = "Neutrophils"
$ws.Cells.Item(30, 5).Value = 3
$ws.Cells.Item(30, 6).Value = 1
$ws.Cells.Item(30, 7).Value = 277.206945
$ws.Cells.Item(30, 8).Value = 831.6208349999999
$ws.Cells.Item(30, 9).Value = 0.551630954676992
$ws.Cells.Item(30, 10).Value = 0.5528566965057786
$ws.Cells.Item(30, 11).Value = 3
$ws.Cells.Item(30, 12).Value = 1
$ws.Cells.Item(30, 13).Value = 52.09024033333333
$ws.Cells.Item(30, 14).Value = 156.270721
$ws.Cells.Item(30, 15).Value = 0.08877815602319267
$ws.Cells.Item(30, 16).Value = 0.09101737745760805
$ws.Cells.Item(30, 17).Value = 14439.77638711911
$ws.Cells.Item(30, 18).Value = 129957.987484072
$ws.Cells.Item(30, 19).Value = 0.04897277896153671
$ws.Cells.Item(30, 20).Value = 0.05031956662583271

# Row 31: Resolving-Mac -> Resolving-Mac
$ws.Cells.Item(31, 1).Value = "Resolving-Mac"
$ws.Cells.Item(31, 2).Value = "C1qb"
$ws.Cells.Item(31, 3).Value = "Lrp1"
$ws.Cells.Item(31, 4).Value = "Resolving-Mac"
$ws.Cells.Item(31, 5).Value = 3
$ws.Cells.Item(31, 6).Value = 1
$ws.Cells.Item(31, 7).Value = 277.206945
$ws.Cells.Item(31, 8).Value = 831.6208349999999
$ws.Cells.Item(31, 9).Value = 0.551630954676992
$ws.Cells.Item(31, 10).Value = 0.5528566965057786
$ws.Cells.Item(31, 11).Value = 3
$ws.Cells.Item(31, 12).Value = 1
$ws.Cells.Item(31, 13).Value = 201.5345866666667
$ws.Cells.Item(31, 14).Value = 604.60376
$ws.Cells.Item(31, 15).Value = 0.3434783342267227
$ws.Cells.Item(31, 16).Value = 0.3521417722019025
$ws.Cells.Item(31, 17).Value = 55866.78708170439
$ws.Cells.Item(31, 18).Value = 502801.0837353395
$ws.Cells.Item(31, 19).Value = 0.18947328142035
$ws.Cells.Item(31, 20).Value = 0.1946839368812343
